$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Columns present in the data rows (B..F, then I..N ; column H/G unused for data writes here
# G holds constant 1 values already present and unchanged; column H is empty)
$cols = @(2,3,4,5,6,9,10,11,12,13,14)  # B,C,D,E,F,I,J,K,L,M,N

$data = @(
    @(1.02,1.033366029227542,1.030263256528264,1.041296853725589,1.049333393264743,1.029555289975002,1.038491151223592,1.033074774661393,1.044076768678686,1.052090743175996,1.016359438818662),
    @(1.02,1.035456500123056,1.030606612076188,1.043210675826196,1.051452051245787,1.029546965970128,1.040218527726439,1.033227698754684,1.045798295184488,1.054018255724803,1.016972109842518),
    @(1.02,1.036804338530668,1.030828240658591,1.044444645201334,1.052818661432289,1.029539753632458,1.041331346126495,1.0333253549493,1.046907440178718,1.055260824523872,1.017365833892181),
    @(1.02,1.037369841693726,1.030921280219693,1.044962378131683,1.053392182731227,1.029536283262567,1.041798024285772,1.0333660975551,1.047372600427122,1.055782111355798,1.017530713907021),
    @(1.02,1.03746472686782,1.030936894093019,1.045049248248012,1.053488421571141,1.029535674837201,1.04187631483863,1.033372920065691,1.047450637651736,1.055869574615349,1.01755836060787),
    @(1.02,1.036811899204402,1.030829484383891,1.0444515671777,1.052826328749198,1.029539708984831,1.0413375863911,1.033325900582437,1.046913660057774,1.055267794227486,1.017368039536576),
    @(1.02,1.034073532951003,1.030379406564799,1.041944565930772,1.050050311958026,1.02955285428455,1.039075957746021,1.033126723317437,1.04465957474345,1.052743136263629,1.016567061135038),
    @(1.02,1.029209801128964,1.029582245551909,1.037492025308604,1.045124365226372,1.029562081003761,1.035051958522594,1.032765901988229,1.040649743577527,1.048257462646569,1.015134458890293),
    @(1.02,1.025939617291363,1.029048220613477,1.034498556500015,1.041815550827719,1.029558927936093,1.032341682945349,1.032518842347416,1.037949528146233,1.045240498997322,1.014164613834184),
    @(1.02,1.024516597594767,1.028816400681362,1.033196024245224,1.040376494231811,1.029555369739793,1.031161200832587,1.032410340936442,1.036773553826061,1.043927448861385,1.013741033680938),
    @(1.02,1.023986934918851,1.028730207336787,1.032711220915079,1.039840979805725,1.02955371971951,1.030721647574989,1.032369811708957,1.036335698255753,1.043438686032643,1.013583141053325),
    @(1.02,1.024100599274606,1.028748699911505,1.032815257965506,1.039955894605354,1.029554088495152,1.030815982143403,1.032378515607504,1.036429667602273,1.043543574808476,1.013617034894785),
    @(1.02,1.024472837959735,1.028809277633772,1.033155970530062,1.040332248738053,1.029555240040407,1.031124889223288,1.032406995397466,1.036737382070019,1.043887068933161,1.013727993633358),
    @(1.02,1.024702040976964,1.028846590368058,1.033365763308621,1.040564001381583,1.029555906069242,1.031315074403884,1.032424512714592,1.036926835544606,1.044098568512557,1.013796284965319),
    @(1.02,1.02603390762909,1.029063593663636,1.034584864946816,1.041910920141611,1.029559117975006,1.032419879352862,1.032526011298701,1.038027428624988,1.04532749814416,1.014192648055195),
    @(1.02,1.02686745018941,1.029199559676678,1.035347855050772,1.042754090270336,1.029560546404202,1.03311102144153,1.032589272239751,1.038715969682559,1.046096560209594,1.014440296640642),
    @(1.02,1.027352967761147,1.029278809869972,1.035792284708298,1.043245289960371,1.029561167926103,1.033513487843186,1.032626024211051,1.039116933547572,1.046544498368184,1.014584396272636),
    @(1.02,1.027518403541067,1.029305822478735,1.035943721312831,1.043412674438442,1.029561343914426,1.033650606598615,1.032638530682611,1.039253542438254,1.04669712568265,1.014633471543508),
    @(1.02,1.026778088910315,1.029184977633817,1.035266056715309,1.042663689154085,1.029560415031907,1.033036937450378,1.032582500137902,1.038642163243583,1.046014113847933,1.014413762527972),
    @(1.02,1.024363253327552,1.028791441332969,1.03305566662043,1.040221449297817,1.029554909994161,1.03103395356437,1.032398615060048,1.036646796958579,1.043785947385488,1.01369533448756),
    @(1.02,1.022838626403462,1.028543518342299,1.031660192195374,1.038680200795731,1.029549549582982,1.029768393883549,1.03228168691861,1.035386161784353,1.042378990933125,1.013240407459075),
    @(1.02,1.023647471882663,1.028674992704237,1.032400512166615,1.039497799546493,1.029552570889732,1.030439889504659,1.032343796520178,1.036055033862706,1.043125427002243,1.013481881958444),
    @(1.02,1.026818469506129,1.029191566806252,1.035303019744479,1.042704539402101,1.029560475047559,1.033070414871011,1.032585560613934,1.038675515205249,1.046051369794693,1.01442575323124),
    @(1.02,1.030471933579147,1.029788797171082,1.03864741647204,1.046402085227728,1.029561340424493,1.03609700809811,1.032860338570589,1.041691023941352,1.049421663745579,1.015507385343153)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $rowIndex = $startRow + $i
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $cols.Count; $j++) {
        $ws.Cells.Item($rowIndex, $cols[$j]).Value2 = $rowVals[$j]
    }
}
